# NIT-9012891781.xlsx — "Elimna EC anteriores y se agregan nuevos, se
# modifica base de datos"
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-62, cols E & F) is
# rebuilt with the period list in reverse-chronological order (newest
# period 2401 first, oldest 2003 last) instead of chronological order.
# Column B/C/D (Tipo Doc, N Doc, Nombre) and column G (Salario Basico)
# are constant across the table and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chronological list of periods as they appeared (row 16 .. row 62).
$periods = @(
    "2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
    "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
    "2401"
)

# Matching "Valor Mora" amounts (row 16 .. row 62) before the edit.
$valores = @(
    35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,
    35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,
    35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,
    35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,35112,
    10534
)

$firstRow = 16
$count = $periods.Length

# Reverse the two columns together (newest period / its value moves to
# the top, oldest moves to the bottom), row position & styling untouched.
for ($i = 0; $i -lt $count; $i++) {
    $row = $firstRow + $i
    $srcIdx = $count - 1 - $i

    $ws.Cells.Item($row, 5).Value = $periods[$srcIdx]
    $ws.Cells.Item($row, 6).Value = $valores[$srcIdx]
}
